$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98, shifting the existing rows 98-100 down to 99-101
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new data record
$ws.Cells.Item(98, 1).Value = 3
$ws.Cells.Item(98, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(98, 3).Value = "Coquimbo"
$ws.Cells.Item(98, 4).Value = 45239
$ws.Cells.Item(98, 5).Value = 5
$ws.Cells.Item(98, 6).Value = 100112022
$ws.Cells.Item(98, 7).Value = "Arveja Verde"
$ws.Cells.Item(98, 8).Value = "Perfection"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 35
$ws.Cells.Item(98, 11).Value = 20000
$ws.Cells.Item(98, 12).Value = 20000
$ws.Cells.Item(98, 13).Value = 20000
$ws.Cells.Item(98, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(98, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(98, 16).Value = 800
$ws.Cells.Item(98, 17).Value = 25
$ws.Cells.Item(98, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date/time number format used by the other rows in column D
$ws.Cells.Item(98, 4).NumberFormat = $ws.Cells.Item(99, 4).NumberFormat()

Write-Host "Row inserted and populated successfully"
